$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A so the existing Date/TransactId/...
# data shifts one column to the right (A:F -> B:G).
$ws.Columns.Item(1).Insert()

# New header cell for the inserted "Bank" column, styled like the other
# header cells (bold header style used by B1:G1).
$ws.Range("A1").Value = "Bank"
$ws.Range("A1").Font.Bold = $true

# Fill the new column with the bank name for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "BankA"
}

# Reset the cursor to A1 (the workbook no longer keeps the old D12
# selection from before the edit).
$ws.Range("A1").Select() | Out-Null
